# EPPlusTest/Workbooks/ComboFromExcel.xlsx - re-save as produced by Excel:
#  - switch workbook calculation back to Automatic (drops calcMode="manual"
#    concurrentManualCount="1" from <calcPr>)
#  - recalculate the volatile RANDBETWEEN() formulas in B20:D42 so their
#    cached <v> results are refreshed
#  - explicitly write "autoZero" c:crosses on both chart axes
#  - move/resize the combo chart (new two-cell anchor)
#  - update the sheet view's selection / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- workbook calculation mode: manual -> automatic -------------------
$excel.Calculation = -4105          # xlCalculationAutomatic

# --- recalc everything (refreshes the RANDBETWEEN cached values) -----
$excel.CalculateFull()

# --- chart tweaks -------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$catAx = $chart.Axes(1)
$valAx = $chart.Axes(2)
$catAx.Crosses = "autoZero"
$valAx.Crosses = "autoZero"

# reposition/resize the chart to its new anchor (from col0/row1 to col8/row16)
$co.Left = 28.87496062992126
$co.Top = 28.87496062992126
$co.Width = 443.5
$co.Height = 216.0

# --- sheet view: move selection, which also drops the stale scroll ----
$ws.Range("K15").Select()
